# Applies the cryptocurrency price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: safe to assign directly via .Value
# (Excel will keep these as text because they cannot be parsed as plain numbers).
$textUpdates = [ordered]@{
    'D2' = '27.476.34'
    'E2' = '  +0.32%  '
    'D3' = '1.635.82'
    'E3' = '  -0.91%  '
    'E4' = '  +0.07%  '
    'E5' = '  -0.44%  '
    'E6' = '  +4.51%  '
    'E7' = '  +0.09%  '
    'E8' = '  -5.16%  '
    'E9' = '  -2.45%  '
    'E10' = '  -0.85%  '
    'E11' = '  +1.42%  '
    'D12' = '1.865.94'
    'E12' = '  -0.96%  '
    'D13' = '1.648.10'
    'E13' = '  +0.21%  '
    'B14' = 'Polygon'
    'C14' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E14' = '  -1.05%  '
    'B15' = 'Polkadot'
    'C15' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'E15' = '  -1.51%  '
    'E16' = '  -2.50%  '
    'D17' = '27.405.91'
    'E17' = '  +0.06%  '
    'E18' = '  -2.25%  '
    'E19' = '  +2.99%  '
    'D20' = '0.0₃0724'
    'E20' = '  -0.28%  '
    'E21' = '  +0.16%  '
    'E22' = '  -2.12%  '
    'E23' = '  +6.64%  '
    'E24' = '  -3.77%  '
    'E25' = '  +1.92%  '
    'E26' = '  -3.28%  '
    'E27' = '  +1.59%  '
    'E28' = '  +0.21%  '
    'E29' = '  -3.16%  '
    'E30' = '  -0.85%  '
    'E31' = '  -1.96%  '
    'E32' = '  -0.58%  '
    'E33' = '  +2.27%  '
    'D34' = '1.407.05'
    'E34' = '  -3.59%  '
    'E35' = '  +2.40%  '
    'E36' = '  -1.82%  '
    'E37' = '  -0.72%  '
    'B38' = 'ARBITRUM'
    'C38' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'E38' = '  -3.63%  '
    'B39' = 'VeChain'
    'C39' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E39' = '  -1.46%  '
    'E40' = '  +13.13%  '
    'E41' = '  -0.33%  '
    'E42' = '  +0.16%  '
    'E44' = '  +1.31%  '
    'E45' = '  +1.35%  '
    'E46' = '  -0.59%  '
    'D47' = '1.775.41'
    'E47' = '  -0.93%  '
    'E48' = '  -3.07%  '
    'E49' = '  -2.68%  '
    'E50' = '  +0.60%  '
    'E51' = '  -1.97%  '
}

foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Values that look like plain numbers (e.g. "212.47", "0.532") must be forced to
# stay text, matching the source data which stores every price as a text string
# (note some prices use "." as a thousands separator, e.g. "27.476.34", and this
# column must remain uniformly textual). We temporarily switch the cell to the
# Text number format before assigning the value, then restore the default "Normal"
# style so no stray formatting is left behind on the cell.
$forcedTextUpdates = [ordered]@{
    'D5' = '212.47'
    'D6' = '0.532'
    'D8' = '22.93'
    'D11' = '0.0888'
    'D14' = '0.564'
    'D15' = '4.02'
    'D16' = '64.24'
    'D18' = '229.05'
    'D19' = '7.67'
    'D22' = '4.31'
    'D23' = '9.91'
    'D24' = '1.94'
    'D25' = '149.74'
    'D26' = '6.96'
    'D29' = '15.56'
    'D32' = '3.29'
    'D33' = '3.17'
    'D37' = '0.570'
    'D38' = '0.876'
    'D39' = '0.0167'
    'D40' = '0.889'
    'D46' = '64.77'
    'D49' = '86.04'
    'D51' = '0.0988'
}

foreach ($ref in $forcedTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $forcedTextUpdates[$ref]
    $cell.Style = "Normal"
}

Write-Host ("Applied {0} text updates and {1} forced-text updates." -f $textUpdates.Count, $forcedTextUpdates.Count)
